$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Provincia de Linares, 2023-09-25) is inserted
# as row 67, pushing the previous rows 67-68 down to 68-69 (their values
# are preserved unchanged by the row insert).
$ws.Rows("67:67").Insert()

$ws.Range("A67").Value = 3
$ws.Range("B67").Value = "Femacal de La Calera"
$ws.Range("C67").Value = "Coquimbo"
$ws.Range("D67").Value = 45194
$ws.Range("E67").Value = 5
$ws.Range("F67").Value = 300000000
$ws.Range("G67").Value = "Espárragos"
$ws.Range("H67").Value = "Verde"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 1300
$ws.Range("K67").Value = 1900
$ws.Range("L67").Value = 1900
$ws.Range("M67").Value = 1900
$ws.Range("N67").Value = "$/kilo"
$ws.Range("O67").Value = "Provincia de Linares"
$ws.Range("P67").Value = 1900
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"
